# Add a "substring in templates" draft example right after the
# ${Root.param1} run that lives in its own single-cell table.
# "${Root.param1}" becomes "${Root.param1} x ${Root.param1[1]}" where
# the "x" is colored red, and every new chunk of text is kept as its
# own run (mirroring how the diff splits formatting boundaries).

$d = $word.ActiveDocument

# Find the 1x1 table whose only cell holds "${Root.param1}" - that is
# the paragraph the diff is anchored to (there's another, identical
# looking "${Root.param1}" paragraph outside any table that must stay
# untouched).
$cell = $null
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $t = $d.Tables.Item($ti)
    if ($t.Rows.Count -eq 1 -and $t.Columns.Count -eq 1) {
        $c = $t.Cell(1, 1)
        if ($c.Range.Text.StartsWith('${Root.param1}')) {
            $cell = $c
        }
    }
}

# Locate the exact text range of "${Root.param1}" inside that cell.
$probe = $d.Range($cell.Range.Start, $cell.Range.Start)
$probe.Find.Execute('${Root.param1}', $false, $false, $false, $false, $false, $true, 0) | Out-Null

$probe.Collapse(0)
$insertStart = $probe.Start

# Insert all of the new text in one go; the host merges identically
# formatted neighbours into a single run, so afterwards we surgically
# re-split the boundaries that must stay separate runs.
$probe.InsertAfter(' x ${Root.param1[1]}')

# --- color the "x" red; this also naturally splits it away from its
#     neighbours on both sides ---
$xRange = $d.Range($insertStart + 1, $insertStart + 2)
$xRange.Font.Color = 255

# --- force additional run boundaries without leaving any formatting
#     residue behind: toggling a boolean property on then back off
#     splits the run at that position, but its final rPr is identical
#     to what it would have been anyway ---
function Split-RunAt($start, $end) {
    $rng = $d.Range($start, $end)
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
}

# "${Root.param1}" / " " boundary (keep the original run separate from
# the new leading space run)
$boundary0Start = $insertStart
$boundary0End = $insertStart + 1
Split-RunAt $boundary0Start $boundary0End

# " ${Root.param1" / "[1]" boundary
$boundary1Start = $insertStart + 3
$boundary1End = $insertStart + 16
Split-RunAt $boundary1Start $boundary1End

# "[1]" / "}" boundary
$boundary2Start = $insertStart + 16
$boundary2End = $insertStart + 19
Split-RunAt $boundary2Start $boundary2End

# "}" final run
$boundary3Start = $insertStart + 19
$boundary3End = $insertStart + 20
Split-RunAt $boundary3Start $boundary3End
